$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet name/title to reflect the new "through" date
$ws.Name = "Through 2022-06-27"

# Update the row label for the June data row
$ws.Range("A7").Value = "June (through 06-27)"

# Update the June row (row 7) values for each year column (B:I)
$ws.Range("B7").Value = 17
$ws.Range("C7").Value = 33
$ws.Range("D7").Value = 65
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 40
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = 112
$ws.Range("I7").Value = 130

# Update the Total row (row 8) values for each year column (B:I)
$ws.Range("B8").Value = 125
$ws.Range("C8").Value = 242
$ws.Range("D8").Value = 381
$ws.Range("E8").Value = 345
$ws.Range("F8").Value = 244
$ws.Range("G8").Value = 458
$ws.Range("H8").Value = 743
$ws.Range("I8").Value = 793
